# RPA datasets push 2024-04-28
# The "하나31호스팩" record (row 11: 하나, 2024-02-22, 하나31호스팩, 하나, 하나,
# 2024-02-27, 2024-03-05, 10000, 5000000, 2000, 0, 100) is removed from the
# table. Deleting the entire row shifts the following row (한국 / 삼현) up to
# become row 11, and the sheet's used range shrinks from A1:L12 to A1:L11.
# The now-unused shared strings for that row are dropped automatically when
# the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Delete()
